# Requirements Trace Matrix - update AutoFilter criteria on the
# "Requirement" column of Table5, and update the saved cell selection.
#
# The table's AutoFilter previously filtered column D ("Status") for the
# value "P". The new filter instead targets column B ("Requirement"),
# restricting the visible rows to four specific EPSG / GeoTIFF related
# requirement statements. Applying the filter through the Excel object
# model also recomputes which data rows are hidden, matching Excel's own
# behaviour when a filter is (re)applied interactively.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)
$rng = $lo.Range

$criteria = @(
    "If the ProjectionMethodGeoKey value is 32767 (User-Defined) then the ProjectedCitationGeoKey and keys for each map projection parameter (coordinate operation parameter) appropriate to that method SHALL be populated.",
    "ProjectionGeoKey values in the range 1024-32766 SHALL be valid EPSG map projection (coordinate operation) codes",
    "ProjMethodGeoKey values in the range 1-27 SHALL be GeoTIFF map projection method codes (coordinate operation method)",
    "The GTModelTypeGeoKey value SHALL be:`n- 0 to indicate that the Model CRS is undefined or unknown`n- 1 to indicate that the Model CRS is a 2D projected Coordinate Reference System, indicated by the value of the ProjectedCRSGeoKey; or`n- 2 to indicate that the Model CRS is a 2DD geographic coordinate reference system, indicated by the value of the GeodeticCRSGeoKey; or`n- 3 to indicate that the Model CRS is a geocentric Cartesian 3D coordinate reference system, indicated by the value of the GeodeticCRSGeoKey; or`n- 32767 to indicate that the Model CRS type is user-defined."
)

# Field 2 = "Requirement" column (1-based field index within the table).
# xlFilterValues = 7, so Criteria1 is treated as an array of discrete
# values to match (equivalent to checking several boxes in the AutoFilter
# dropdown).
$rng.AutoFilter(2, $criteria, 7) | Out-Null

# Move the saved selection to B71, matching the author's last position.
$ws.Range("B71").Select() | Out-Null
